$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The last data row (row 28, unit code E_EWPKM2) needs to move up to become
# row 5 (right after the E_BNEUR row), pushing the existing rows 5-27 down
# by one row each (to rows 6-28).
#
# Note: we use Value2 (not Value) to read/write cell contents - in this
# COM model, reading the Value property directly returns the property's
# reflection signature instead of invoking it, while Value2 behaves like
# the classic Excel COM property.

$firstRow = 5
$lastRow = 28

# Read current contents of rows 5..28 (columns A, B, C) into arrays first,
# so that later writes don't clobber values we still need to read.
$valsA = @()
$valsB = @()
$valsC = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $valsA += $ws.Range("A$r").Value2
    $valsB += $ws.Range("B$r").Value2
    $valsC += $ws.Range("C$r").Value2
}

$count = $lastRow - $firstRow + 1
$lastIdx = $count - 1

# The row that used to be last (index $lastIdx, originally row 28) becomes
# the new row 5.
$ws.Range("A$firstRow").Value2 = $valsA[$lastIdx]
$ws.Range("B$firstRow").Value2 = $valsB[$lastIdx]
$ws.Range("C$firstRow").Value2 = $valsC[$lastIdx]

# All the other rows (originally rows 5..27, i.e. array indices 0..$lastIdx-1)
# shift down by one row, to rows 6..28.
for ($i = 0; $i -lt $lastIdx; $i++) {
    $rn = $firstRow + 1 + $i
    $ws.Range("A$rn").Value2 = $valsA[$i]
    $ws.Range("B$rn").Value2 = $valsB[$i]
    $ws.Range("C$rn").Value2 = $valsC[$i]
}
